{"js": "// Update the date line and the 25 multiplication answers in the table.\n// Each old value is unique within the document, so a scoped search +\n// in-place replace (which preserves the run's formatting) is safe.\nconst replacements = [\n  [\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"],\n  [\"76\u00d787=6612\", \"11\u00d763=693\"],\n  [\"54\u00d756=3024\", \"73\u00d788=6424\"],\n  [\"35\u00d784=2940\", \"58\u00d741=2378\"],\n  [\"44\u00d773=3212\", \"52\u00d714=728\"],\n  [\"25\u00d789=2225\", \"37\u00d737=1369\"],\n  [\"55\u00d712=660\", \"17\u00d729=493\"],\n  [\"43\u00d745=1935\", \"66\u00d768=4488\"],\n  [\"47\u00d796=4512\", \"69\u00d793=6417\"],\n  [\"80\u00d732=2560\", \"79\u00d792=7268\"],\n  [\"85\u00d798=8330\", \"58\u00d717=986\"],\n  [\"26\u00d778=2028\", \"41\u00d744=1804\"],\n  [\"62\u00d712=744\", \"91\u00d769=6279\"],\n  [\"79\u00d763=4977\", \"28\u00d793=2604\"],\n  [\"91\u00d712=1092\", \"40\u00d772=2880\"],\n  [\"37\u00d756=2072\", \"36\u00d754=1944\"],\n  [\"34\u00d735=1190\", \"75\u00d766=4950\"],\n  [\"29\u00d719=551\", \"36\u00d762=2232\"],\n  [\"80\u00d715=1200\", \"98\u00d716=1568\"],\n  [\"14\u00d796=1344\", \"81\u00d754=4374\"],\n  [\"75\u00d728=2100\", \"20\u00d761=1220\"],\n  [\"32\u00d756=1792\", \"88\u00d794=8272\"],\n  [\"90\u00d713=1170\", \"67\u00d752=3484\"],\n  [\"32\u00d764=2048\", \"97\u00d712=1164\"],\n  [\"81\u00d799=8019\", \"38\u00d752=1976\"],\n  [\"34\u00d721=714\", \"79\u00d778=6162\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication answers in the table.\n# Each old value is unique within the document, so Find/Replace (scoped to\n# the whole document body) safely targets exactly one run each time and\n# preserves that run's formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-12 Tuesday\", \"2025-08-13 Wednesday\"),\n    @(\"76\u00d787=6612\", \"11\u00d763=693\"),\n    @(\"54\u00d756=3024\", \"73\u00d788=6424\"),\n    @(\"35\u00d784=2940\", \"58\u00d741=2378\"),\n    @(\"44\u00d773=3212\", \"52\u00d714=728\"),\n    @(\"25\u00d789=2225\", \"37\u00d737=1369\"),\n    @(\"55\u00d712=660\", \"17\u00d729=493\"),\n    @(\"43\u00d745=1935\", \"66\u00d768=4488\"),\n    @(\"47\u00d796=4512\", \"69\u00d793=6417\"),\n    @(\"80\u00d732=2560\", \"79\u00d792=7268\"),\n    @(\"85\u00d798=8330\", \"58\u00d717=986\"),\n    @(\"26\u00d778=2028\", \"41\u00d744=1804\"),\n    @(\"62\u00d712=744\", \"91\u00d769=6279\"),\n    @(\"79\u00d763=4977\", \"28\u00d793=2604\"),\n    @(\"91\u00d712=1092\", \"40\u00d772=2880\"),\n    @(\"37\u00d756=2072\", \"36\u00d754=1944\"),\n    @(\"34\u00d735=1190\", \"75\u00d766=4950\"),\n    @(\"29\u00d719=551\", \"36\u00d762=2232\"),\n    @(\"80\u00d715=1200\", \"98\u00d716=1568\"),\n    @(\"14\u00d796=1344\", \"81\u00d754=4374\"),\n    @(\"75\u00d728=2100\", \"20\u00d761=1220\"),\n    @(\"32\u00d756=1792\", \"88\u00d794=8272\"),\n    @(\"90\u00d713=1170\", \"67\u00d752=3484\"),\n    @(\"32\u00d764=2048\", \"97\u00d712=1164\"),\n    @(\"81\u00d799=8019\", \"38\u00d752=1976\"),\n    @(\"34\u00d721=714\", \"79\u00d778=6162\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
